$wb = $excel.ActiveWorkbook

$wsIndexClose = $wb.Worksheets.Item("Index Close")
$wsMTD = $wb.Worksheets.Item("MTD %")
$wsDoD = $wb.Worksheets.Item("Day over Day %")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsDailyMovers = $wb.Worksheets.Item("Daily Movers")
$wsMarketOverview = $wb.Worksheets.Item("Market Overview")

# --- Sheet: Index Close ---
$wsIndexClose.Range("G2").Value = 60310.1484375
$wsIndexClose.Range("B7").Value = 25898.55078125
$wsIndexClose.Range("C7").Value = 67939.5
$wsIndexClose.Range("D7").Value = 26402.44921875
$wsIndexClose.Range("F7").Value = 17089.75
$wsIndexClose.Range("G7").Value = 59578.05078125
$wsIndexClose.Range("I7").Value = 23550.849609375
$wsIndexClose.Range("K7").Value = 23550.849609375

# --- Sheet: MTD % ---
$wsMTD.Range("G2").Value = 6.13
$wsMTD.Range("B7").Value = 5.17
$wsMTD.Range("C7").Value = 1.9
$wsMTD.Range("D7").Value = 4.64
$wsMTD.Range("F7").Value = 6.59
$wsMTD.Range("G7").Value = 4.84
$wsMTD.Range("I7").Value = 3.69
$wsMTD.Range("K7").Value = 3.69

# --- Sheet: Day over Day % ---
$wsDoD.Range("G2").Value = 1.2
$wsDoD.Range("B6").Value = 0.5
$wsDoD.Range("C6").Value = 1.01
$wsDoD.Range("D6").Value = 0.58
$wsDoD.Range("F6").Value = 1.01
$wsDoD.Range("G6").Value = 1.07
$wsDoD.Range("I6").Value = 0.74
$wsDoD.Range("K6").Value = 0.74
$wsDoD.Range("B7").Value = 0
$wsDoD.Range("C7").Value = 0
$wsDoD.Range("D7").Value = 0
$wsDoD.Range("F7").Value = 0
$wsDoD.Range("G7").Value = 0
$wsDoD.Range("I7").Value = 0
$wsDoD.Range("K7").Value = 0

# --- Sheet: Summary (rows shifted) ---
$wsSummary.Range("A3").Value = "Nifty Midcap 100"
$wsSummary.Range("B3").Value = 6.13
$wsSummary.Range("A4").Value = "Nifty 50"
$wsSummary.Range("B4").Value = 5.45
$wsSummary.Range("A5").Value = "Nifty 200"
$wsSummary.Range("B5").Value = 5.2
$wsSummary.Range("A6").Value = "Nifty 100"
$wsSummary.Range("B6").Value = 4.98
$wsSummary.Range("A7").Value = "Nifty Midcap 150"
$wsSummary.Range("B7").Value = 4.9

# --- Sheet: Daily Movers ---
$wsDailyMovers.Range("B2").Value = "Nifty Midcap 100, Nifty Midcap 150, Nifty Midcap 50"
$wsDailyMovers.Range("C2").Value = "Nifty 50, Nifty 100, Nifty 200"
$wsDailyMovers.Range("B6").Value = "Nifty Midcap 100, Nifty Next 50, Nifty Midcap 50"
$wsDailyMovers.Range("C6").Value = "Nifty Midcap 150, Nifty 200, Nifty500 Multicap 50:25:25"
$wsDailyMovers.Range("B7").Value = "Nifty Midcap 150, Nifty500 Multicap 50:25:25, Nifty 200"
$wsDailyMovers.Range("C7").Value = "Nifty 50, Nifty Next 50, Nifty 100"

# --- Sheet: Market Overview ---
# Force text so Excel doesn't auto-convert the percent string into a numeric
# percentage value/format, then clear the temporary number format so the
# cell's style matches its unstyled neighbours.
$wsMarketOverview.Range("B3").NumberFormat = "@"
$wsMarketOverview.Range("B3").Value = "4.89%"
$wsMarketOverview.Range("B3").ClearFormats()
